$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 465
$ws.Range("I2").Value = 1232
$ws.Range("J2").Value = 4852
$ws.Range("K2").Value = 26
$ws.Range("L2").Value = 1331
$ws.Range("M2").Value = 95
$ws.Range("N2").Value = 845
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 18
$ws.Range("Q2").Value = 15
$ws.Range("R2").Value = 48
$ws.Range("S2").Value = 520
$ws.Range("T2").Value = 846
$ws.Range("U2").Value = 65
$ws.Range("V2").Value = 7672
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 7731
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 129
$ws.Range("AA2").Value = 36
